$wb = $excel.ActiveWorkbook

# The same F2/F3/F4 updates apply to both the "展览" and "全部类型" sheets,
# which hold identical data (want-to-go counts for each event).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1564
    $ws.Range("F3").Value = 89
    $ws.Range("F4").Value = 25
}
